$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("AREA_POR_REGIAO")
$ws.Range("O14").Value = 0
